$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# 1. ". Although XRootD supported ..." -> "The XRootD ability to federate ..."
ReplaceText ". Although XRootD supported multi-storage deployments for a long time, the addition of a feature that allowed its proper functionality within a global, multi-site environment was in fact the core idea of AAA." ". The XRootD ability to federate different sides through meta managers together with additional functionalities provided by the AAA (like logical file name translation to physical file name) allowed to achieve a global, multi-site environment for data storage and analysis."

Write-Output "step1 done"

# 2. Big "A user might want to retrieve ..." TCP/epoll paragraph rewrite
$oldTcp = "A user might want to retrieve some data using the XRootD client from a file that is located on a server. Interaction between the XRootD client and that particular server is done through a TCP implementation. The mechanism that allows the client to receive feedback from the TCP kernel is called event-loop. The feedback consists of communicating whether there is available space in the TCP-output buffer for writing data (i.e., requests which will be sent to the server) or if there is some data in the TCP-receive buffer for reading responses from the server. In this event-based workflow, there is a queue of requests that the client is issuing to the server, and with each write-event, a request is removed from the queue and it is being written on the socket. It is worth mentioning that the TCP buffers (for both sending and receiving data) might not have enough size to allow requests/responses to be written/read in a single event, meaning that it can take several write/read events to process an entire request to the server or a response from the server. Furthermore, each request has a corresponding message handler, so that after a request is written to the socket (in order to be sent to the server), the accompanying message handler is moved into a queue for incoming responses. During a read-event yielded by the event-loop, the client is informed that it can read from the socket, that is a server response. Once a response arrived from the server, its corresponding message handler (located inside the incoming queue) is also taken out from the queue, and finally, after the response is parsed, the callback function is being called."
$newTcp = "A user might want to retrieve some data using the XRootD client from a file that is located on a server. Interaction between the XRootD client and that particular server is done over TCP protocol. Using the epoll system call, the XRootD client runtime receives events from the kernel signaling whether there is available space in the TCP-output buffer for writing data (i.e., requests which will be sent to the server) or if there is some data in the TCP-receive buffer for reading responses from the server. In this event-based workflow, there is a queue of requests that the client is issuing to the server, and with each write-event, a request is removed from the queue and it is being written on the socket. It is worth mentioning that the TCP buffers (for both sending and receiving data) might not have enough size to allow requests/responses to be written/read in a single event, meaning that it can take several write/read events to process an entire request to the server or a response from the server. Furthermore, each request has a corresponding message handler, so that after a request is written to the socket (in order to be sent to the server), the accompanying message handler is moved into a queue for incoming responses. During a read-event yielded by the event-loop, the client is informed that it can readout from the socket the server response. Once a response arrived from the server, its corresponding message handler (located inside the incoming queue) is also taken out from the queue, and finally, after the response is parsed, the callback function is being called."
ReplaceText $oldTcp $newTcp

# Italicize the newly-inserted "epoll" word
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("epoll", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $rng2.Italic = 1
}

Write-Output "step2 done"

# 3. "It is in fact the response handler ..." -> "In case the user wishes to use only asynchronous operations, ..."
ReplaceText " It is in fact the response handler that takes care of the function callback once it has been executed; in other words, the handler controls the proper flow of the execution pipeline. The flow of operations follows works in such a way that each next function from the pipeline needs to be called within the handler of the previous function." " In case the user wishes to use only asynchronous operations, the subsequent operation needs to be called from the handler of the previous operation."

Write-Output "step3 done"

# 4. "e.g." -> "e.g.," (both occurrences) and "must have" -> "must call"
ReplaceText "might consist of a function that tries to open the file (e.g. Open). Its response handler must have the second operation (e.g. Read, Write) that needs to be called." "might consist of a function that tries to open the file (e.g., Open). Its response handler must call the second operation (e.g., Read, Write) that needs to be called."

Write-Output "step4 done"

# 5. "The constructed API makes it so there is a communication protocol between the operations: " -> "The proposed API provides a syntax for chaining consecutive operations: "
ReplaceText "The constructed API makes it so there is a communication protocol between the operations: " "The proposed API provides a syntax for chaining consecutive operations: "

Write-Output "step5 done"

# 6. "The defined operations are connected to" -> "The defined operations are chained to"
ReplaceText "The defined operations are connected to each other by the | operator." "The defined operations are chained to each other by the | operator."

Write-Output "step6 done"

# 7. "In order to emphasize the overall flexibility and fluidity of the pipeline syntax" -> remove "and fluidity"
ReplaceText "In order to emphasize the overall flexibility and fluidity of the pipeline syntax, the following example is proposed: " "In order to emphasize the overall flexibility of the pipeline syntax, the following example is proposed: "

Write-Output "step7 done"

# 8. remove "(taking as an argument the lock file itself)"
ReplaceText "declaration of a lock file, then the lock file is created with the first call of the Open function (taking as an argument the lock file itself). Once the lock file has been created," "declaration of a lock file, then the lock file is created with the first call of the Open function. Once the lock file has been created,"

Write-Output "step8 done"

# 9. "The Rm function is used" -> "The Rm() function is used"
ReplaceText "file that needs to be accessed. The Rm function is used for deleting the lock file since it is not needed anymore." "file that needs to be accessed. The Rm() function is used for deleting the lock file since it is not needed anymore."

Write-Output "step9 done"

# 10. "The Declarative API is tested in the development of an Erasure Coding plug-in for the client." -> "The main use case for the declarative API is the development of an erasure coding plugin for the client."
ReplaceText "The Declarative API is tested in the development of an Erasure Coding plug-in for the client. Erasure Coding (EC) is a method of data protection" "The main use case for the declarative API is the development of an erasure coding plugin for the client. Erasure Coding (EC) is a method of data protection"

Write-Output "step10 done"

# 11. Big "goal of erasure coding" paragraph rewrite
$oldEc = "The goal of erasure coding is to enable data that becomes corrupted at some point in the disk storage process to be reconstructed by using information about the data that's stored elsewhere in the array. The drawback of erasure coding is that it can be more CPU-intensive, and that can translate into increased latency. In other words, erasure coding adds the redundancy to the system that tolerates failures. In terms of the workflow, EC takes the original data and encodes it in such a way that when needed, only a subset of all the chunks is required to recreate the original information. The data protection scheme is graphically represented in Figs. 5 and 6, where the decode and encode procedures, respectively, are explained."
$newEc = "The goal of erasure coding is to enable data that becomes corrupted at some point in the disk storage process to be reconstructed by using information about the data that's stored elsewhere in the array. The tradeoff for erasure coding is that it can be more CPU-intensive, and that error recovery might result in increased network traffic and latency. EC encodes N chunks of data (of equal size) in such a way that the result is the N original data chunks and additional K chunks of parity (N+K chunks in total). Every N chunks of the obtained N+K chunks are sufficient to recover the original N chunks. The data protection scheme is graphically represented in Figs. 5 and 6, where the decode and encode procedures, respectively, are explained."
ReplaceText $oldEc $newEc

Write-Output "step11 done"

# 12. "asynchronous operations. The process of writing the plug-in achieves..." paragraph
$oldAsync = "asynchronous operations. The process of writing the plug-in achieves a high degree of code readability, with a clear workflow and reduced complexity. The standard asynchronous operations hide the actual workflow of operations behind the first function callback"
$newAsync = "asynchronous operations. The obtained code is much more readable, with a clear workflow and reduced complexity. In the contrary, the standard asynchronous operations hide the actual workflow of operations behind the first function callback"
ReplaceText $oldAsync $newAsync

Write-Output "step12 done"
